$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 500 and 501),
# pushing the existing rows 500-613 down to become 502-615.
$ws.Rows.Item(500).Insert()
$ws.Rows.Item(500).Insert()

# --- New row 500 ---
$ws.Range("A500").Value = 9
$ws.Range("B500").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C500").Value = "Metropolitana"
$ws.Range("D500").Value = 44798
$ws.Range("D500").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E500").Value = 13
$ws.Range("F500").Value = 100114014
$ws.Range("G500").Value = "Betarraga"
$ws.Range("H500").Value = "Sin especificar"
$ws.Range("I500").Value = "Primera"
$ws.Range("J500").Value = 4300
$ws.Range("K500").Value = 160
$ws.Range("L500").Value = 170
$ws.Range("M500").Value = 165
$ws.Range("N500").Value = "`$/unidad"
$ws.Range("O500").Value = "Región Metropolitana"
$ws.Range("P500").Value = 165
$ws.Range("Q500").Value = 1
$ws.Range("R500").Value = "Hortaliza"

# --- New row 501 ---
$ws.Range("A501").Value = 9
$ws.Range("B501").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C501").Value = "Metropolitana"
$ws.Range("D501").Value = 44798
$ws.Range("D501").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E501").Value = 13
$ws.Range("F501").Value = 100114014
$ws.Range("G501").Value = "Betarraga"
$ws.Range("H501").Value = "Sin especificar"
$ws.Range("I501").Value = "Segunda"
$ws.Range("J501").Value = 9700
$ws.Range("K501").Value = 140
$ws.Range("L501").Value = 140
$ws.Range("M501").Value = 140
$ws.Range("N501").Value = "`$/unidad"
$ws.Range("O501").Value = "Región Metropolitana"
$ws.Range("P501").Value = 140
$ws.Range("Q501").Value = 1
$ws.Range("R501").Value = "Hortaliza"
